# Insert two new price-record rows into the daily "Zapallo italiano"
# (Feria Lagunitas de Puerto Montt) sheet:
#   - one new row before the (old) row 324, becoming the new row 324
#   - one new row before the (old) row 424, becoming the new row 424
#     (after the first insertion it sits right before the row that used
#     to be row 424 and has by then been pushed down to row 425)
# Every other existing row keeps its data and simply shifts down by the
# number of new rows inserted above it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert first new row, at position 324 -------------------------------
$ws.Rows.Item(324).Insert()

# --- Insert second new row, at position 424 (post first-insert numbering)
$ws.Rows.Item(424).Insert()

# --- Populate the new row 324 --------------------------------------------
$ws.Cells.Item(324, 1).Value = 4
$ws.Cells.Item(324, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(324, 3).Value = "Los Lagos"
$ws.Cells.Item(324, 4).Value = 45120
$ws.Cells.Item(324, 5).Value = 10
$ws.Cells.Item(324, 6).Value = 100112032
$ws.Cells.Item(324, 7).Value = "Zapallo italiano"
$ws.Cells.Item(324, 8).Value = "Sin especificar"
$ws.Cells.Item(324, 9).Value = "Primera"
$ws.Cells.Item(324, 10).Value = 100
$ws.Cells.Item(324, 11).Value = 20000
$ws.Cells.Item(324, 12).Value = 22000
$ws.Cells.Item(324, 13).Value = 21000
$ws.Cells.Item(324, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(324, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(324, 16).Value = 420
$ws.Cells.Item(324, 17).Value = 50
$ws.Cells.Item(324, 18).Value = "Hortaliza"

# --- Populate the new row 424 --------------------------------------------
$ws.Cells.Item(424, 1).Value = 4
$ws.Cells.Item(424, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(424, 3).Value = "Los Lagos"
$ws.Cells.Item(424, 4).Value = 45121
$ws.Cells.Item(424, 5).Value = 10
$ws.Cells.Item(424, 6).Value = 100112032
$ws.Cells.Item(424, 7).Value = "Zapallo italiano"
$ws.Cells.Item(424, 8).Value = "Sin especificar"
$ws.Cells.Item(424, 9).Value = "Primera"
$ws.Cells.Item(424, 10).Value = 240
$ws.Cells.Item(424, 11).Value = 19000
$ws.Cells.Item(424, 12).Value = 20000
$ws.Cells.Item(424, 13).Value = 19500
$ws.Cells.Item(424, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(424, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(424, 16).Value = 390
$ws.Cells.Item(424, 17).Value = 50
$ws.Cells.Item(424, 18).Value = "Hortaliza"

# --- Fix up the D column's date number format / style on the new rows ----
$ws.Cells.Item(324, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(424, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
